# Apply "data added to excel table" edits to the "conferences" sheet:
#  - ATBC2023 row (row 2) gains a code-of-conduct link, a revised expected
#    attendants figure, a revised keynote-speaker count, a new female
#    keynote-speaker count, a contact e-mail and a comment.
#  - Two new one-cell rows are appended for ECONET2023 and PHENOLOGY2022.
#  - A brand-new ATBC2022 conference row is appended with full details.
#  - The "website code of conduct" header is renamed "website_codeconduct".

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("conferences")

# --- header rename -------------------------------------------------------
$ws.Range("L1").Value = "website_codeconduct"

# --- ATBC2023 (row 2) updates / additions ---------------------------------
$ws.Range("L2").Value = "https://www.atbc2023.org/codeofconduct"
$ws.Range("M2").Value = ">400"
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = "atbc2023@tropicalbio.org"
$ws.Range("Q2").Value = "Wrote them asking for children facilities"

# --- new single-cell conference rows ---------------------------------------
$ws.Range("A7").Value = "ECONET2023"
$ws.Range("A8").Value = "PHENOLOGY2022"

# --- new ATBC2022 conference row (row 9) ------------------------------------
$ws.Range("A9").Value = "ATBC2022"
$ws.Range("B9").Value = "https://www.atbc2022.org/"
$ws.Range("C9").Value = "Cartagena"
$ws.Range("D9").Value = "Colombia"

# dates: copy the existing date-format from row 2 so the new cells keep the
# same number format / style as the other start-date / end-date cells
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E9").Value = 44752

$ws.Range("F2").Copy() | Out-Null
$ws.Range("F9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F9").Value = 44756

$ws.Range("G9").Value = "yes"

# organisation name: same wrapped text/style as the ATBC2023 row above
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H9").Value = $ws.Range("H2").Value()

$ws.Range("I9").Value = $ws.Range("I2").Value()
$ws.Range("J9").Value = "Conserving Tropical Biodiversity and Achieving Socio-Ecological Resilience in the Anthropocene:" + [char]0x200B + " Opportunities and Challenges"
$ws.Range("K9").Value = "yes"
$ws.Range("L9").Value = "https://www.atbc2022.org/codeofconduct"

$excel.CutCopyMode = $false

# leave the selection where the author's last edit was
$ws.Activate() | Out-Null
$ws.Range("L9").Select() | Out-Null
